$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2373
$ws.Range("J51").Value = 2330.6667
$ws.Range("L51").Value = 2330.6667
$ws.Range("N51").Value = -3298.6667
$ws.Range("H64").Value = 3421.923
$ws.Range("J64").Value = 3622
$ws.Range("L64").Value = 3622
$ws.Range("N64").Value = -4118
$ws.Range("H67").Value = 3421.923
$ws.Range("J67").Value = 3622
$ws.Range("L67").Value = 3622
$ws.Range("N67").Value = -5338
$ws.Range("H116").Value = 6389.5557
$ws.Range("I116").Value = 2500
$ws.Range("J116").Value = 6875.75
$ws.Range("K116").Value = 2500
$ws.Range("L116").Value = 6875.75
$ws.Range("M116").Value = 942
$ws.Range("N116").Value = -13759.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 19999
$ws.Range("J76").Value = 19999
$ws.Range("L76").Value = 19999
$ws.Range("N76").Value = -20675
$ws.Range("H79").Value = 19999
$ws.Range("J79").Value = 19999
$ws.Range("L79").Value = 19999
$ws.Range("N79").Value = -22339
$ws.Range("H92").Value = 15552.25
$ws.Range("J92").Value = 15552.25
$ws.Range("L92").Value = 15552.25
$ws.Range("N92").Value = -20544.25
$ws.Range("H132").Value = 15024.815
$ws.Range("I132").Value = 1390.7587
$ws.Range("J132").Value = 58956.777
$ws.Range("K132").Value = 4172.2761
$ws.Range("L132").Value = 176870.331
$ws.Range("M132").Value = -1642.2761
$ws.Range("N132").Value = -181930.331

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1824.3334
$ws.Range("I94").Value = 952.25
$ws.Range("J94").Value = 2522
$ws.Range("K94").Value = 952.25
$ws.Range("L94").Value = 2522
$ws.Range("M94").Value = -501.25
$ws.Range("N94").Value = -3424
$ws.Range("H99").Value = 1448.7273
$ws.Range("I99").Value = 1437.3334
$ws.Range("K99").Value = 1437.3334
$ws.Range("M99").Value = 60.66660000000002

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 66671276
$ws.Range("I62").Value = 100003910
$ws.Range("J62").Value = 6002.4
$ws.Range("K62").Value = 100003910
$ws.Range("L62").Value = 6002.4
$ws.Range("M62").Value = -100003286
$ws.Range("N62").Value = -7250.4
$ws.Range("H65").Value = 66671276
$ws.Range("I65").Value = 100003910
$ws.Range("J65").Value = 6002.4
$ws.Range("K65").Value = 500019550
$ws.Range("L65").Value = 30012
$ws.Range("M65").Value = -500016430
$ws.Range("N65").Value = -36252
$ws.Range("H105").Value = 785.84
$ws.Range("I105").Value = 717.0454999999999
$ws.Range("K105").Value = 717.0454999999999
$ws.Range("M105").Value = 1029.9545

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 58.933334
$ws.Range("I12").Value = 4.8
$ws.Range("J12").Value = 86
$ws.Range("K12").Value = 14.4
$ws.Range("L12").Value = 258
$ws.Range("M12").Value = 158.6
$ws.Range("N12").Value = -604
$ws.Range("H23").Value = 606.8333
$ws.Range("I23").Value = 20.5
$ws.Range("J23").Value = 900
$ws.Range("K23").Value = 61.5
$ws.Range("L23").Value = 2700
$ws.Range("M23").Value = 173.5
$ws.Range("N23").Value = -3170
$ws.Range("H33").Value = 68.166664
$ws.Range("I33").Value = 20
$ws.Range("J33").Value = 77.8
$ws.Range("K33").Value = 120
$ws.Range("L33").Value = 466.8
$ws.Range("M33").Value = 163
$ws.Range("N33").Value = -1032.8
$ws.Range("H68").Value = 1736.6154
$ws.Range("J68").Value = 1847.6
$ws.Range("L68").Value = 5542.799999999999
$ws.Range("N68").Value = -7164.799999999999
$ws.Range("H71").Value = 1736.6154
$ws.Range("J71").Value = 1847.6
$ws.Range("L71").Value = 16628.4
$ws.Range("N71").Value = -24740.4
$ws.Range("H82").Value = 8211.6
$ws.Range("I82").Value = 1013
$ws.Range("J82").Value = 10011.25
$ws.Range("K82").Value = 3039
$ws.Range("L82").Value = 30033.75
$ws.Range("M82").Value = -2633
$ws.Range("N82").Value = -30845.75
$ws.Range("H85").Value = 8211.6
$ws.Range("I85").Value = 1013
$ws.Range("J85").Value = 10011.25
$ws.Range("K85").Value = 3039
$ws.Range("L85").Value = 30033.75
$ws.Range("M85").Value = -1635
$ws.Range("N85").Value = -32841.75
$ws.Range("H96").Value = 8205.714
$ws.Range("J96").Value = 9323.333000000001
$ws.Range("L96").Value = 27969.999
$ws.Range("N96").Value = -32087.999
$ws.Range("H131").Value = 749.99
$ws.Range("I131").Value = 625
$ws.Range("J131").Value = 752.54083
$ws.Range("K131").Value = 1875
$ws.Range("L131").Value = 2257.62249
$ws.Range("M131").Value = 3165
$ws.Range("N131").Value = -12337.62249
$ws.Range("H133").Value = 4883
$ws.Range("I133").Value = 2119.75
$ws.Range("J133").Value = 6111.1113
$ws.Range("K133").Value = 6359.25
$ws.Range("L133").Value = 18333.3339
$ws.Range("M133").Value = -1299.25
$ws.Range("N133").Value = -28453.3339

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1058.1666
$ws.Range("I22").Value = 1744.7142
$ws.Range("J22").Value = 621.2727
$ws.Range("K22").Value = 1744.7142
$ws.Range("L22").Value = 621.2727
$ws.Range("M22").Value = -1449.7142
$ws.Range("N22").Value = -1211.2727
$ws.Range("H27").Value = 1058.1666
$ws.Range("I27").Value = 1744.7142
$ws.Range("J27").Value = 621.2727
$ws.Range("K27").Value = 1744.7142
$ws.Range("L27").Value = 621.2727
$ws.Range("M27").Value = -1637.7142
$ws.Range("N27").Value = -835.2727
$ws.Range("H61").Value = 4221.385
$ws.Range("I61").Value = 1807.091
$ws.Range("J61").Value = 17500
$ws.Range("K61").Value = 1807.091
$ws.Range("L61").Value = 17500
$ws.Range("M61").Value = -1605.091
$ws.Range("N61").Value = -17904
$ws.Range("H113").Value = 4221.385
$ws.Range("I113").Value = 1807.091
$ws.Range("J113").Value = 17500
$ws.Range("K113").Value = 1807.091
$ws.Range("L113").Value = 17500
$ws.Range("M113").Value = 362.9090000000001
$ws.Range("N113").Value = -21840
$ws.Range("H132").Value = 1207767
$ws.Range("I132").Value = 3013798
$ws.Range("J132").Value = 3746.3333
$ws.Range("K132").Value = 9041394
$ws.Range("L132").Value = 11238.9999
$ws.Range("M132").Value = -9038864
$ws.Range("N132").Value = -16298.9999

Write-Host "All changes applied successfully."